$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the columns which were swapped when the codeforIATI SectorGroup codelist
# was generated: column E (codeforiati:category-name) and column F
# (codeforiati:group-code) hold each other's data. Swap the header labels and
# every data row back into place, except for the category rows whose
# category-code is itself the "lead" category of its group (111, 121, 151,
# 231, 311, 321) - those rows already carry the correct values.

$excludedCategoryCodes = @("111", "121", "151", "231", "311", "321")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Header row: swap E1 / F1 labels.
$headerE = $ws.Cells.Item(1, 5).Value2
$headerF = $ws.Cells.Item(1, 6).Value2
$ws.Cells.Item(1, 5).Value2 = $headerF
$ws.Cells.Item(1, 6).Value2 = $headerE

for ($r = 2; $r -le $lastRow; $r++) {
    $categoryCode = $ws.Cells.Item($r, 4).Value2
    if ($excludedCategoryCodes -notcontains $categoryCode) {
        $eVal = $ws.Cells.Item($r, 5).Value2
        $fVal = $ws.Cells.Item($r, 6).Value2
        $ws.Cells.Item($r, 5).Value2 = $fVal
        $ws.Cells.Item($r, 6).Value2 = $eVal
    }
}
